$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.820.05'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '1.629.33'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5071'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06430'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07810'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.255'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').Value = '1.627.40'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').Value = '1.853.76'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5576'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.27'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('D17').Value = '0.0₅7519'
$ws.Range('E17').Value = '  -2.89%  '
$ws.Range('D18').Value = '25.817.99'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.298'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.25%  '
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('E25').Value = '  -3.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1275'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '140.29'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.718'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.238'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.282'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.185'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.553'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8936'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.28%  '
$ws.Range('D37').Value = '1.135.95'
$ws.Range('E37').Value = '  +3.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.550'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5457'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.67%  '
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.000'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.571'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7946'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.22'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.98%  '
$ws.Range('D45').Value = '1.778.39'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').Value = '0.0₈112'
$ws.Range('E46').Value = '  -7.56%  '
$ws.Range('E47').Value = '  -2.10%  '
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('E49').Value = '  -2.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.580'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('D51').Style = 'Normal'
